$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1721.6897
$ws.Range("J17").Value = 1721.6897
$ws.Range("L17").Value = 5165.0691
$ws.Range("N17").Value = -5501.0691
$ws.Range("H19").Value = 2060.3157
$ws.Range("J19").Value = 2208.6667
$ws.Range("L19").Value = 2208.6667
$ws.Range("N19").Value = -2558.6667
$ws.Range("H80").Value = 1328.6666
$ws.Range("I80").Value = 1243.7142
$ws.Range("J80").Value = 1447.6
$ws.Range("K80").Value = 3731.1426
$ws.Range("L80").Value = 4342.799999999999
$ws.Range("M80").Value = -2733.1426
$ws.Range("N80").Value = -6338.799999999999
$ws.Range("H83").Value = 1328.6666
$ws.Range("I83").Value = 1243.7142
$ws.Range("J83").Value = 1447.6
$ws.Range("K83").Value = 11193.4278
$ws.Range("L83").Value = 13028.4
$ws.Range("M83").Value = -6201.427799999999
$ws.Range("N83").Value = -23012.4
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H106").Value = 32978.9
$ws.Range("I106").Value = 36309.89
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 36309.89
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -35678.89
$ws.Range("N106").Value = -4262
$ws.Range("H112").Value = 4186.6875
$ws.Range("J112").Value = 4186.6875
$ws.Range("L112").Value = 12560.0625
$ws.Range("N112").Value = -14776.0625
$ws.Range("H113").Value = 144890.42
$ws.Range("I113").Value = 1322
$ws.Range("J113").Value = 288458.84
$ws.Range("K113").Value = 1322
$ws.Range("L113").Value = 288458.84
$ws.Range("M113").Value = 1932
$ws.Range("N113").Value = -294966.84
$ws.Range("H116").Value = 3800.3333
$ws.Range("J116").Value = 4469
$ws.Range("L116").Value = 4469
$ws.Range("N116").Value = -11353
$ws.Range("H137").Value = 4726165.5
$ws.Range("I137").Value = 8932097
$ws.Range("K137").Value = 26796291
$ws.Range("M137").Value = -26793741

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 1117.1111
$ws.Range("J50").Value = 1272.2858
$ws.Range("L50").Value = 1272.2858
$ws.Range("N50").Value = -2700.2858
$ws.Range("H74").Value = 260914.83
$ws.Range("I74").Value = 325395.47
$ws.Range("J74").Value = 11052.25
$ws.Range("K74").Value = 325395.47
$ws.Range("L74").Value = 11052.25
$ws.Range("M74").Value = -324521.47
$ws.Range("N74").Value = -12800.25
$ws.Range("H77").Value = 260914.83
$ws.Range("I77").Value = 325395.47
$ws.Range("J77").Value = 11052.25
$ws.Range("K77").Value = 1626977.35
$ws.Range("L77").Value = 55261.25
$ws.Range("M77").Value = -1622609.35
$ws.Range("N77").Value = -63997.25
$ws.Range("H104").Value = 60056
$ws.Range("J104").Value = 60056
$ws.Range("L104").Value = 60056
$ws.Range("N104").Value = -67044

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 774.75
$ws.Range("J7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("N7").Value = -1226
$ws.Range("H20").Value = 4145.875
$ws.Range("I20").Value = 5502
$ws.Range("K20").Value = 5502
$ws.Range("M20").Value = -5255
$ws.Range("H94").Value = 702.1429000000001
$ws.Range("I94").Value = 653.5714
$ws.Range("K94").Value = 653.5714
$ws.Range("M94").Value = -202.5714
$ws.Range("H105").Value = 333418200
$ws.Range("I105").Value = 500125000
$ws.Range("J105").Value = 4600
$ws.Range("K105").Value = 500125000
$ws.Range("L105").Value = 4600
$ws.Range("M105").Value = -500123253
$ws.Range("N105").Value = -8094
$ws.Range("H134").Value = 7013.5947
$ws.Range("I134").Value = 6291.9565
$ws.Range("J134").Value = 8199.143
$ws.Range("K134").Value = 18875.8695
$ws.Range("L134").Value = 24597.429
$ws.Range("M134").Value = -16340.8695
$ws.Range("N134").Value = -29667.429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 289.8
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 149.66667
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 149.66667
$ws.Range("M12").Value = -330
$ws.Range("N12").Value = -489.66667
$ws.Range("H31").Value = 32260964
$ws.Range("I31").Value = 50000870
$ws.Range("K31").Value = 50000870
$ws.Range("M31").Value = -50000575
$ws.Range("H34").Value = 32260964
$ws.Range("I34").Value = 50000870
$ws.Range("K34").Value = 50000870
$ws.Range("M34").Value = -50000668
$ws.Range("H63").Value = 69989
$ws.Range("J63").Value = 69989
$ws.Range("L63").Value = 69989
$ws.Range("N63").Value = -71361
$ws.Range("H66").Value = 69989
$ws.Range("J66").Value = 69989
$ws.Range("L66").Value = 209967
$ws.Range("N66").Value = -216831
$ws.Range("H140").Value = 83439
$ws.Range("J140").Value = 83439
$ws.Range("L140").Value = 83439
$ws.Range("N140").Value = -93799

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2672.3
$ws.Range("I114").Value = 2121.1667
$ws.Range("J114").Value = 3499
$ws.Range("K114").Value = 6363.500100000001
$ws.Range("L114").Value = 10497
$ws.Range("M114").Value = -3109.500100000001
$ws.Range("N114").Value = -17005
$ws.Range("H129").Value = 26319352
$ws.Range("I129").Value = 4620.5713
$ws.Range("J129").Value = 100000600
$ws.Range("K129").Value = 13861.7139
$ws.Range("L129").Value = 300001800
$ws.Range("M129").Value = -8861.713899999999
$ws.Range("N129").Value = -300011800

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 42498
$ws.Range("J48").Value = 42498
$ws.Range("L48").Value = 42498
$ws.Range("N48").Value = -43468
$ws.Range("H70").Value = 37640.145
$ws.Range("I70").Value = 35830.168
$ws.Range("K70").Value = 35830.168
$ws.Range("M70").Value = -35560.168
$ws.Range("H73").Value = 37640.145
$ws.Range("I73").Value = 35830.168
$ws.Range("K73").Value = 35830.168
$ws.Range("M73").Value = -34894.168
$ws.Range("H92").Value = 34600
$ws.Range("J92").Value = 34600
$ws.Range("L92").Value = 34600
$ws.Range("N92").Value = -38344
$ws.Range("H97").Value = 1385.7778
$ws.Range("I97").Value = 1440.8
$ws.Range("K97").Value = 1440.8
$ws.Range("M97").Value = -944.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1671.125
$ws.Range("I16").Value = 1731.2858
$ws.Range("J16").Value = 1250
$ws.Range("K16").Value = 1731.2858
$ws.Range("L16").Value = 1250
$ws.Range("M16").Value = -1561.2858
$ws.Range("N16").Value = -1590
$ws.Range("H22").Value = 3272.5588
$ws.Range("I22").Value = 1767.7059
$ws.Range("J22").Value = 4777.4116
$ws.Range("K22").Value = 1767.7059
$ws.Range("L22").Value = 4777.4116
$ws.Range("M22").Value = -1472.7059
$ws.Range("N22").Value = -5367.4116
$ws.Range("H27").Value = 3272.5588
$ws.Range("I27").Value = 1767.7059
$ws.Range("J27").Value = 4777.4116
$ws.Range("K27").Value = 1767.7059
$ws.Range("L27").Value = 4777.4116
$ws.Range("M27").Value = -1660.7059
$ws.Range("N27").Value = -4991.4116
$ws.Range("H93").Value = 1153.0667
$ws.Range("I93").Value = 1054
$ws.Range("K93").Value = 1054
$ws.Range("M93").Value = 194
$ws.Range("H94").Value = 59165
$ws.Range("J94").Value = 59165
$ws.Range("L94").Value = 59165
$ws.Range("N94").Value = -60517

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 10044
$ws.Range("I39").Value = 10044
$ws.Range("K39").Value = 10044
$ws.Range("M39").Value = -9631
$ws.Range("H122").Value = 2086.9375
$ws.Range("I122").Value = 1811.1
$ws.Range("K122").Value = 5433.299999999999
$ws.Range("M122").Value = -2983.299999999999
$ws.Range("H136").Value = 4307.263
$ws.Range("I136").Value = 1126.6364
$ws.Range("J136").Value = 8680.625
$ws.Range("K136").Value = 3379.9092
$ws.Range("L136").Value = 26041.875
$ws.Range("M136").Value = -829.9092000000001
$ws.Range("N136").Value = -31141.875
